$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where DAMSLTag changed from "sd" (Statement-non-opinion) to "sv" (Statement-opinion)
$svRows = @(29, 31, 37, 41)
foreach ($r in $svRows) {
    $ws.Range("I$r").Value = "sv"
    $ws.Range("J$r").Value = "Statement-opinion"
}

# Rows where DAMSLTag changed to "aa" (Agree/Accept)
$ws.Range("I38").Value = "aa"
$ws.Range("J38").Value = "Agree/Accept"

$ws.Range("I44").Value = "aa"
$ws.Range("J44").Value = "Agree/Accept"

$wb.Save()
